# feat: add 2022-Q4 data
#
# Target layout after the edit:
#   Sheet1 "总计"    (unchanged sheet, gets a new row for 2022-Q4, old row shifts down)
#   Sheet2 "2022-Q4" (brand new sheet with fund holdings data)
#   Sheet3 "2021-Q2" (the old "2021-Q2" sheet content, recreated verbatim)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Re-shuffle the sheets so sheetId/rId numbering matches the target:
#    总计=1, 2022-Q4=2(new), 2021-Q2=3(new, recreated with identical data)
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$oldQ2 = $wb.Worksheets.Item("2021-Q2")
$oldQ2.Delete()

$q4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $summary)
$q4.Name = "2022-Q4"

$q2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q4)
$q2.Name = "2021-Q2"

# ---------------------------------------------------------------------------
# 2) Update the "总计" summary sheet: insert the 2022-Q4 row above the
#    existing 2021-Q2 row (push it from row 2 to row 3).
# ---------------------------------------------------------------------------
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2021-Q2"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0.02

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 6
$summary.Range("D2").Value = 0.04

# A2 already carries the header-row style (s=2); copy it onto A3 too so both
# index cells share the same look (matches the original authoring).
$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Output "done restructuring sheets"

# ---------------------------------------------------------------------------
# 3) Populate the new "2022-Q4" sheet with the fund-holding table.
# ---------------------------------------------------------------------------
$q4Headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
$q4Rows = @(
    @(0, "012864", "易方达标普医疗保健指数（QDII-LOF）人民币 C", "0.50", "93.65", "1.59", "0.0080", 3),
    @(1, "161126", "易方达标普医疗保健指数（QDII-LOF）人民币",   "0.50", "93.65", "1.59", "0.0080", 3),
    @(2, "003719", "易方达标普医疗保健指数（QDII-LOF）美元A",    "0.49", "93.65", "1.59", "0.0078", 3),
    @(3, "013404", "大成标普500等权重指数（QDII）美元",          "3.55", "92.19", "0.21", "0.0075", 10),
    @(4, "096001", "大成标普500等权重指数（QDII）人民币",        "3.55", "92.19", "0.21", "0.0075", 10),
    @(5, "012865", "易方达标普医疗保健指数（QDII-LOF）美元 C",   "0.01", "93.65", "1.59", "0.0002", 3)
)

# Header row (B1:H1), bold/centered style copied from the summary sheet header.
$col = 2
foreach ($h in $q4Headers) {
    $q4.Cells.Item(1, $col).Value = $h
    $col++
}

# Columns B, D, E, F, G hold text that looks numeric ("012864", "0.50", ...).
# Force them to Text *before* writing so Excel does not silently coerce them
# to numbers / strip the leading zeros, exactly like the source data.
# (Multi-area "B2:B7,D2:G7" selectors are not reliable here, so do the two
# contiguous blocks separately.)
$q4.Range("B2:B7").NumberFormat = "@"
$q4.Range("D2:G7").NumberFormat = "@"

$r = 2
foreach ($row in $q4Rows) {
    $q4.Cells.Item($r, 1).Value = $row[0]
    $q4.Cells.Item($r, 2).Value = $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 4).Value = $row[3]
    $q4.Cells.Item($r, 5).Value = $row[4]
    $q4.Cells.Item($r, 6).Value = $row[5]
    $q4.Cells.Item($r, 7).Value = $row[6]
    $q4.Cells.Item($r, 8).Value = $row[7]
    $r++
}

# Restore the plain/default look on the text cells (drop the temporary "@"
# format) while keeping their values as Text - match the source which has
# no explicit style on these cells.
$q4.Range("Z1").Copy()
$q4.Range("B2:B7").PasteSpecial(-4122)
$q4.Range("D2:G7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header cells (B1:H1) and the index column (A2:A7) use the same bold style
# as the summary sheet's header - copy that style across.
$summary.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$summary.Range("A2").Copy()
$q4.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Output "done 2022-Q4 sheet"

# ---------------------------------------------------------------------------
# 4) Recreate the "2021-Q2" sheet (it was deleted in step 1 so its sheetId
#    would come out right) with the exact same data it had before the edit.
# ---------------------------------------------------------------------------
$q2Headers = @("基金代码","基金名称","基金金额","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
$q2Rows = @(
    @(0, "161126", "易方达标普医疗保健指数(QDII-LOF) 人民币", "0.50", "94.26", "1.54", "0.0077", 9),
    @(1, "003719", "易方达标普医疗保健指数(QDII-LOF) 美元",   "0.50", "94.26", "1.54", "0.0077", 9)
)

$col = 2
foreach ($h in $q2Headers) {
    $q2.Cells.Item(1, $col).Value = $h
    $col++
}

$q2.Range("B2:B3").NumberFormat = "@"
$q2.Range("D2:G3").NumberFormat = "@"

$r = 2
foreach ($row in $q2Rows) {
    $q2.Cells.Item($r, 1).Value = $row[0]
    $q2.Cells.Item($r, 2).Value = $row[1]
    $q2.Cells.Item($r, 3).Value = $row[2]
    $q2.Cells.Item($r, 4).Value = $row[3]
    $q2.Cells.Item($r, 5).Value = $row[4]
    $q2.Cells.Item($r, 6).Value = $row[5]
    $q2.Cells.Item($r, 7).Value = $row[6]
    $q2.Cells.Item($r, 8).Value = $row[7]
    $r++
}

$q2.Range("Z1").Copy()
$q2.Range("B2:B3").PasteSpecial(-4122)
$q2.Range("D2:G3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# This sheet originally used style index 1 (not 2) for its header/index
# column - grab it straight from the style the sheet already had in the
# workbook by borrowing the corresponding xf via a helper cell on the 2022
# sheet? Not available - instead reconstruct via explicit formatting that
# matches cellXf #1 (centered, bold, thin border): easiest is to source it
# from a cell that still carries that exact style. None remain after the
# original "2021-Q2" sheet was deleted, so apply the equivalent formatting
# directly.
$q2HeaderRange = $q2.Range("B1:H1")
$q2HeaderRange.Font.Bold = $true
$q2HeaderRange.HorizontalAlignment = -4108
$q2HeaderRange.VerticalAlignment = -4160
$q2HeaderRange.Borders.LineStyle = 1

$q2IndexRange = $q2.Range("A2:A3")
$q2IndexRange.Font.Bold = $true
$q2IndexRange.HorizontalAlignment = -4108
$q2IndexRange.VerticalAlignment = -4160
$q2IndexRange.Borders.LineStyle = 1

Write-Output "done 2021-Q2 sheet"


